$d = $word.ActiveDocument

# Locate the "url" meta-tag line by its distinctive old URL text.
$oldUrlText = " https://starsindust.github.io/ I-Snuck-A-Book/PDF_Optimizer.html"
$fullText = $d.Content.Text
$startPos = $fullText.IndexOf($oldUrlText)
if ($startPos -lt 0) {
    throw "Could not locate target URL text in document"
}

# Find the paragraph that contains this text so we know exactly where its
# visible content ends (Paragraph.Range.End includes the trailing pilcrim,
# so the last real character is at End - 1).
$paraEnd = -1
foreach ($p in $d.Paragraphs) {
    $pr = $p.Range
    if ($startPos -ge $pr.Start -and $startPos -lt $pr.End) {
        $paraEnd = $pr.End - 1
        break
    }
}
if ($paraEnd -lt 0) {
    throw "Could not locate containing paragraph"
}

# Replace everything from the start of the old URL run through the end of
# the paragraph's content (i.e. " https://...PDF_Optimizer.html"&gt;) with
# the new run layout:
#   " "  "https://starsindust.github.io"  "/"
#   " Enlightenment/Articles/2025/4_Game_Maker_2/3_Creating_Objects/3_Creating_Objects.html "
#   """&gt;"
$target = $d.Range($startPos, $paraEnd)

$payload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
              '<w:r><w:t>https://starsindust.github.io</w:t></w:r>' +
              '<w:r><w:t>/</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> Enlightenment/Articles/2025/4_Game_Maker_2/3_Creating_Objects/3_Creating_Objects.html </w:t></w:r>' +
              '<w:r><w:t>"&gt;</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$target.InsertXML($payload)
